# Apply updated bus voltage magnitude (vm_pu) results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.0516334435754
$ws.Range("D2").Value = 1.056359688791404
$ws.Range("E2").Value = 1.064506053812791
$ws.Range("F2").Value = 1.070263906301469
$ws.Range("I2").Value = 1.045761696423113
$ws.Range("J2").Value = 1.05666028622786
$ws.Range("K2").Value = 1.059097615050719
$ws.Range("L2").Value = 1.067221807516844
$ws.Range("M2").Value = 1.072964208639068
$ws.Range("N2").Value = 1.058160865034933

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052547439563288
$ws.Range("D3").Value = 1.057070996826647
$ws.Range("E3").Value = 1.065379058889307
$ws.Range("F3").Value = 1.071152103033569
$ws.Range("I3").Value = 1.045981689177238
$ws.Range("J3").Value = 1.057224616786036
$ws.Range("K3").Value = 1.059622807162075
$ws.Range("L3").Value = 1.067909886069966
$ws.Range("M3").Value = 1.073668555944095
$ws.Range("N3").Value = 1.058725997007231

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053139521911875
$ws.Range("D4").Value = 1.057531831537789
$ws.Range("E4").Value = 1.065944952713894
$ws.Range("F4").Value = 1.071727860446451
$ws.Range("I4").Value = 1.046123177485899
$ws.Range("J4").Value = 1.057589781268904
$ws.Range("K4").Value = 1.059962536819187
$ws.Range("L4").Value = 1.068355476031912
$ws.Range("M4").Value = 1.074124708551529
$ws.Range("N4").Value = 1.059091680065549

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053388591378931
$ws.Range("D5").Value = 1.05772570153048
$ws.Range("E5").Value = 1.066183092503432
$ws.Range("F5").Value = 1.071970154883403
$ws.Range("I5").Value = 1.046182452242802
$ws.Range("J5").Value = 1.057743296425214
$ws.Range("K5").Value = 1.060105332917648
$ws.Range("L5").Value = 1.068542886175538
$ws.Range("M5").Value = 1.074316567720925
$ws.Range("N5").Value = 1.059245413230989

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053430420467187
$ws.Range("D6").Value = 1.05775826103049
$ws.Range("E6").Value = 1.066223091131823
$ws.Range("F6").Value = 1.072010851562976
$ws.Range("I6").Value = 1.046192392582622
$ws.Range("J6").Value = 1.057769072262237
$ws.Range("K6").Value = 1.060129307437647
$ws.Range("L6").Value = 1.068574358056457
$ws.Range("M6").Value = 1.074348787109195
$ws.Range("N6").Value = 1.059271225672658

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053142849370946
$ws.Range("D7").Value = 1.057534421509332
$ws.Range("E7").Value = 1.065948133816211
$ws.Range("F7").Value = 1.071731097032008
$ws.Range("I7").Value = 1.046123970331972
$ws.Range("J7").Value = 1.057591832547355
$ws.Range("K7").Value = 1.059964444972249
$ws.Range("L7").Value = 1.068357979886593
$ws.Range("M7").Value = 1.074127271821016
$ws.Range("N7").Value = 1.059093734257051

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.051942194120122
$ws.Range("D8").Value = 1.056599959523002
$ws.Range("E8").Value = 1.064800881672661
$ws.Range("F8").Value = 1.070563861275976
$ws.Range("I8").Value = 1.045836222047185
$ws.Range("J8").Value = 1.056851002674816
$ws.Range("K8").Value = 1.059275127170408
$ws.Range("L8").Value = 1.067454272113521
$ws.Range("M8").Value = 1.073202163919148
$ws.Range("N8").Value = 1.058351852321118

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.049831645520124
$ws.Range("D9").Value = 1.054957760077436
$ws.Range("E9").Value = 1.062787017749172
$ws.Range("F9").Value = 1.068515042656207
$ws.Range("I9").Value = 1.045322601123249
$ws.Range("J9").Value = 1.055545654438238
$ws.Range("K9").Value = 1.058059710585825
$ws.Range("L9").Value = 1.065864619488555
$ws.Range("M9").Value = 1.071575077021784
$ws.Range("N9").Value = 1.057044650340319

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048428162162654
$ws.Range("D10").Value = 1.053866041159258
$ws.Range("E10").Value = 1.061449744834115
$ws.Range("F10").Value = 1.067154643266408
$ws.Range("I10").Value = 1.044975806413504
$ws.Range("J10").Value = 1.054675550810682
$ws.Range("K10").Value = 1.057248997849131
$ws.Range("L10").Value = 1.064806812578921
$ws.Range("M10").Value = 1.070492502217861
$ws.Range("N10").Value = 1.05617331106581

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047821296869824
$ws.Range("D11").Value = 1.053394065345949
$ws.Range("E11").Value = 1.060871969399656
$ws.Range("F11").Value = 1.06656689625261
$ws.Range("I11").Value = 1.044824609309895
$ws.Range("J11").Value = 1.05429883126851
$ws.Range("K11").Value = 1.056897861717064
$ws.Range("L11").Value = 1.06434925289549
$ws.Range("M11").Value = 1.070024263074463
$ws.Range("N11").Value = 1.055796056538677

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.047596009442493
$ws.Range("D12").Value = 1.05321886614026
$ws.Range("E12").Value = 1.060657550729353
$ws.Range("F12").Value = 1.066348780098732
$ws.Range("I12").Value = 1.044768293437566
$ws.Range("J12").Value = 1.054158908176455
$ws.Range("K12").Value = 1.056767421439966
$ws.Range("L12").Value = 1.064179367993609
$ws.Range("M12").Value = 1.069850418190377
$ws.Range("N12").Value = 1.055655934739782

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.047644328485202
$ws.Range("D13").Value = 1.053256441797167
$ws.Range("E13").Value = 1.060703535508333
$ws.Range("F13").Value = 1.06639555770994
$ws.Range("I13").Value = 1.044780380366505
$ws.Range("J13").Value = 1.054188921804995
$ws.Range("K13").Value = 1.056795401877437
$ws.Range("L13").Value = 1.064215805550111
$ws.Range("M13").Value = 1.069887704865906
$ws.Range("N13").Value = 1.055685990991117

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047802671894176
$ws.Range("D14").Value = 1.053379580997786
$ws.Range("E14").Value = 1.060854241541741
$ws.Range("F14").Value = 1.066548862628633
$ws.Range("I14").Value = 1.044819957371233
$ws.Range("J14").Value = 1.054287265026105
$ws.Range("K14").Value = 1.056887079735066
$ws.Range("L14").Value = 1.064335208656592
$ws.Range("M14").Value = 1.070009891362261
$ws.Range("N14").Value = 1.055784473870882

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.047900249675897
$ws.Range("D15").Value = 1.053455466227684
$ws.Range("E15").Value = 1.060947122107331
$ws.Range("F15").Value = 1.066643345303269
$ws.Range("I15").Value = 1.044844321611018
$ws.Range("J15").Value = 1.054347858534658
$ws.Range("K15").Value = 1.05694356384308
$ws.Range("L15").Value = 1.064408786583724
$ws.Range("M15").Value = 1.070085185136354
$ws.Range("N15").Value = 1.055845153429167

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048468455882238
$ws.Range("D16").Value = 1.053897380486459
$ws.Range("E16").Value = 1.061488116840219
$ws.Range("F16").Value = 1.067193677949167
$ws.Range("I16").Value = 1.04498581916194
$ws.Range("J16").Value = 1.054700553409335
$ws.Range("K16").Value = 1.057272299755221
$ws.Range("L16").Value = 1.06483718948943
$ws.Range("M16").Value = 1.07052358884589
$ws.Range("N16").Value = 1.05619834917102

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048825105728842
$ws.Range("D17").Value = 1.05417478239267
$ws.Range("E17").Value = 1.061827810418107
$ws.Range("F17").Value = 1.067539240319269
$ws.Range("I17").Value = 1.045074300786911
$ws.Range("J17").Value = 1.054921801392308
$ws.Range("K17").Value = 1.057478483253075
$ws.Range("L17").Value = 1.06510604417136
$ws.Range("M17").Value = 1.070798728907557
$ws.Range("N17").Value = 1.056419911351507

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049033215692462
$ws.Range("D18").Value = 1.054336658071836
$ws.Range("E18").Value = 1.062026070607113
$ws.Range("F18").Value = 1.067740927757397
$ws.Range("I18").Value = 1.045125810939146
$ws.Range("J18").Value = 1.055050855483084
$ws.Range("K18").Value = 1.057598737628091
$ws.Range("L18").Value = 1.065262908556738
$ws.Range("M18").Value = 1.070959263806308
$ws.Range("N18").Value = 1.056549148713894

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049104189704657
$ws.Range("D19").Value = 1.054391865632041
$ws.Range("E19").Value = 1.062093692964809
$ws.Range("F19").Value = 1.067809719409492
$ws.Range("I19").Value = 1.045143357643303
$ws.Range("J19").Value = 1.055094860217109
$ws.Range("K19").Value = 1.057639739705312
$ws.Range("L19").Value = 1.06531640301883
$ws.Range("M19").Value = 1.071014010556116
$ws.Range("N19").Value = 1.056593215939689

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048786832069714
$ws.Range("D20").Value = 1.054145012351114
$ws.Range("E20").Value = 1.061791351811458
$ws.Range("F20").Value = 1.067502151624225
$ws.Range("I20").Value = 1.045064817852127
$ws.Range("J20").Value = 1.054898063167409
$ws.Range("K20").Value = 1.057456362634559
$ws.Range("L20").Value = 1.065077193840715
$ws.Range("M20").Value = 1.070769203764024
$ws.Range("N20").Value = 1.056396139415605

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047756040125154
$ws.Range("D21").Value = 1.053343316422918
$ws.Range("E21").Value = 1.060809857034141
$ws.Range("F21").Value = 1.066503712636927
$ws.Range("I21").Value = 1.044808307190758
$ws.Range("J21").Value = 1.054258305187342
$ws.Range("K21").Value = 1.056860083231773
$ws.Range("L21").Value = 1.064300045387817
$ws.Range("M21").Value = 1.069973908265162
$ws.Range("N21").Value = 1.055755472905826

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047108689430356
$ws.Range("D22").Value = 1.052839916569039
$ws.Range("E22").Value = 1.060193868916171
$ws.Range("F22").Value = 1.065877108100778
$ws.Range("I22").Value = 1.044646134966572
$ws.Range("J22").Value = 1.053856106871153
$ws.Range("K22").Value = 1.056485105900566
$ws.Range("L22").Value = 1.06381184514056
$ws.Range("M22").Value = 1.069474337877846
$ws.Range("N22").Value = 1.055352703421895

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047451790832352
$ws.Range("D23").Value = 1.05310671543409
$ws.Range("E23").Value = 1.060520309475374
$ws.Range("F23").Value = 1.066209173094652
$ws.Range("I23").Value = 1.044732190045486
$ws.Range("J23").Value = 1.054069315356076
$ws.Range("K23").Value = 1.056683894990623
$ws.Range("L23").Value = 1.064070608696507
$ws.Range("M23").Value = 1.069739125185996
$ws.Range("N23").Value = 1.055566214687322

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048804126035132
$ws.Range("D24").Value = 1.054158463929675
$ws.Range("E24").Value = 1.06180782550762
$ws.Range("F24").Value = 1.067518910017604
$ws.Range("I24").Value = 1.045069103090238
$ws.Range("J24").Value = 1.054908789437563
$ws.Range("K24").Value = 1.057466358017747
$ws.Range("L24").Value = 1.065090229921316
$ws.Range("M24").Value = 1.070782544750065
$ws.Range("N24").Value = 1.056406880918294

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050376653270359
$ws.Range("D25").Value = 1.055381771896686
$ws.Range("E25").Value = 1.063306722365009
$ws.Range("F25").Value = 1.069043752947314
$ws.Range("I25").Value = 1.045456159350275
$ws.Range("J25").Value = 1.05588310064609
$ws.Range("K25").Value = 1.058374006683161
$ws.Range("L25").Value = 1.066275243054135
$ws.Range("M25").Value = 1.071995345281163
$ws.Range("N25").Value = 1.057382575760492
